$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the first occurrence of $searchText (from the start of the
# document each time) and replace it with $newText, preserving the
# formatting of the matched run by inserting the new text immediately after
# the match and then deleting the original matched characters.
# ---------------------------------------------------------------------------
function Replace-FirstMatch($doc, $searchText, $newText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $searchText"
        return $null
    }
    $oldLen = $rng.End - $rng.Start
    $startPos = $rng.Start
    $rng.InsertAfter($newText)
    $oldRng = $doc.Range($startPos, $startPos + $oldLen)
    $oldRng.Delete()
    return $startPos
}

# 1) "Welcome to your Python Bootcamp group" + ". " + "You'll ..."
#    (insert ". " right after "group", before "You'll")
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("group", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found1) {
    $rng1.InsertAfter(". ")
}

# 2) Found a data file -> give it its real name, keeping the bold run.
$null = Replace-FirstMatch $d "data file name" "boco_air_temp.csv"

# 3) Fix the "Pands" typo -> "Pandas" (insert the missing "a").
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Pands", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found3) {
    $stem = $d.Range($rng3.Start, $rng3.Start + 4)
    $stem.InsertAfter("a")
}

# 4) Expand the plotting instructions after "plotted against time)".
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("plotted against time)", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found4) {
    $insertion = " of the average annual mean. Also create a plot with the lines for both March and October."
    $rng4.InsertAfter($insertion)
    # The original text continued with ". Once ..." - remove the now-duplicate
    # period since our insertion already ends the sentence with a period.
    $oldDot = $d.Range($rng4.End, $rng4.End + 1)
    if ($oldDot.Text -eq ".") {
        $oldDot.Delete()
    }
}

# 5) "change the line color" -> "change the size of the graph"
$null = Replace-FirstMatch $d "line color" "size of the graph"
